$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in row 1, columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-37
$data = @{
    2  = @(8, 9)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(5, 5)
    8  = @(9, 10)
    9  = @(8, 9)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(8, 9)
    21 = @(8, 8)
    22 = @(7, 8)
    23 = @(7, 8)
    24 = @(6, 8)
    25 = @(3, 6)
    26 = @(3, 5)
    27 = @(1, 3)
    28 = @(6, 8)
    29 = @(5, 7)
    30 = @(7, 8)
    31 = @(1, 5)
    32 = @(1, 6)
    33 = @(1, 4)
    34 = @(1, 6)
    35 = @(1, 4)
    36 = @(1, 3)
    37 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
